$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster data (player, position, team) for rows 2..19 (row 1 is the header).
$data = @(
    @("Jalen Brunson",      "PG",       "New York Knicks"),
    @("Devin Booker",       "PG,SG",    "Phoenix Suns"),
    @("Norman Powell",      "SG,SF",    "LA Clippers"),
    @("Jared McCain",       "PG,SG",    "Philadelphia 76ers"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("P.J. Washington",    "PF",       "Dallas Mavericks"),
    @("Desmond Bane",       "SG,SF",    "Memphis Grizzlies"),
    @("Shaedon Sharpe",     "SG,SF",    "Portland Trail Blazers"),
    @("Goga Bitadze",       "C",        "Orlando Magic"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Walker Kessler",     "C",        "Utah Jazz"),
    @("Dereck Lively II",   "C",        "Dallas Mavericks"),
    @("LeBron James",       "SF,PF",    "Los Angeles Lakers"),
    @("Nicolas Claxton",    "C",        "Brooklyn Nets"),
    @("Trae Young",         "PG",       "Atlanta Hawks"),
    @("Immanuel Quickley",  "PG,SG",    "Toronto Raptors"),
    @("Kawhi Leonard",      "SG,SF,PF", "LA Clippers"),
    @("Coby White",         "PG,SG",    "Chicago Bulls")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
